# Add the Inventory app-service URL next to the "Inventory URL" label (row 17)
# and leave the selection sitting on the new cell, mirroring what a user does
# in Excel after typing the value and pressing Enter/clicking away.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = "my-inventory-abhinav.azurewebsites.net"

[void]$ws.Activate()
[void]$ws.Range("B17").Select()
